# DESIGN/rules/DesignFirst/Main.xlsx - "Rules" sheet
# Commit: "Project DesignFirst is saved." (admin, SAVE)
#
# Semantic change: cell D10 is updated from 21 to 100 (numeric value).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("D10").Value = 100.0
